$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the login rows with freshly generated test credentials (row 5 / row 6)
$ws.Range("C5").Value = "gkqsvvoniujdguz@gmail.com"
$ws.Range("D5").Value = "locnaIKOUE5"

$ws.Range("C6").Value = "covxqsvvpyvvwlx@gmail.com"
$ws.Range("D6").Value = "wrbzmVYMMN5"

# Drop the trailing "pass/PASS" flag on row 6 and remove the now-unused row 7 entirely
$ws.Range("F6").ClearContents() | Out-Null
$ws.Range("C7:E7").ClearContents() | Out-Null

# Widen the email column so the addresses are fully visible
$ws.Columns("C").ColumnWidth = 29.6640625

# Restore a sensible selection/zoom for the cleaned-up sheet
$excel.ActiveWindow.Zoom = 160
$ws.Range("C8").Select() | Out-Null
